$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gal"
$ws.Range("C2").Value = "Galr1"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.776285
$ws.Range("H2").Value = 5.328855000000001
$ws.Range("I2").Value = 0.2250740306326953
$ws.Range("J2").Value = 0.2250740306326953
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03918366666666667
$ws.Range("N2").Value = 0.117551
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.06960135934500002
$ws.Range("R2").Value = 0.6264122341050001
$ws.Range("S2").Value = 0.2250740306326953
$ws.Range("T2").Value = 0.2250740306326953

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gal"
$ws.Range("C3").Value = "Galr1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.407624
$ws.Range("H3").Value = 4.222872
$ws.Range("I3").Value = 0.1783607964348723
$ws.Range("J3").Value = 0.1783607964348722
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.03918366666666667
$ws.Range("N3").Value = 0.117551
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.055155869608
$ws.Range("R3").Value = 0.496402826472
$ws.Range("S3").Value = 0.1783607964348723
$ws.Range("T3").Value = 0.1783607964348722

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Gal"
$ws.Range("C4").Value = "Galr1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.248107333333333
$ws.Range("H4").Value = 12.744322
$ws.Range("I4").Value = 0.5382799720054182
$ws.Range("J4").Value = 0.5382799720054181
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03918366666666667
$ws.Range("N4").Value = 0.117551
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.1664564217135555
$ws.Range("R4").Value = 1.498107795422
$ws.Range("S4").Value = 0.5382799720054182
$ws.Range("T4").Value = 0.5382799720054181

# Row 5
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Gal"
$ws.Range("C5").Value = "Galr1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.459987
$ws.Range("H5").Value = 1.379961
$ws.Range("I5").Value = 0.05828520092701431
$ws.Range("J5").Value = 0.0582852009270143
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03918366666666667
$ws.Range("N5").Value = 0.117551
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.018023977279
$ws.Range("R5").Value = 0.162215795511
$ws.Range("S5").Value = 0.05828520092701431
$ws.Range("T5").Value = 0.0582852009270143

# Remove now-obsolete rows 6-9 (shift cells up)
$ws.Range("A6:T9").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp) | Out-Null

Write-Output "applied edits"